$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 355.5
$ws.Range("I6").Value = 408.25
$ws.Range("J6").Value = 250
$ws.Range("K6").Value = 1224.75
$ws.Range("L6").Value = 750
$ws.Range("M6").Value = -1112.75
$ws.Range("N6").Value = -974

$ws.Range("H33").Value = 464.34616
$ws.Range("I33").Value = 705.73334
$ws.Range("J33").Value = 135.18182
$ws.Range("K33").Value = 705.73334
$ws.Range("L33").Value = 135.18182
$ws.Range("M33").Value = -476.73334
$ws.Range("N33").Value = -593.18182

$ws.Range("H76").Value = 66745800
$ws.Range("I76").Value = 180165
$ws.Range("K76").Value = 180165
$ws.Range("M76").Value = -179850

$ws.Range("H79").Value = 66745800
$ws.Range("I79").Value = 180165
$ws.Range("K79").Value = 180165
$ws.Range("M79").Value = -179073

$ws.Range("H98").Value = 2701.0688
$ws.Range("I98").Value = 796
$ws.Range("K98").Value = 796
$ws.Range("M98").Value = 702

$ws.Range("H122").Value = 2701.0688
$ws.Range("I122").Value = 796
$ws.Range("K122").Value = 2388
$ws.Range("M122").Value = 62

$ws.Range("H132").Value = 5236.5625
$ws.Range("I132").Value = 5388.75
$ws.Range("K132").Value = 16166.25
$ws.Range("M132").Value = -13636.25

$ws.Range("H137").Value = 3678.8262
$ws.Range("I137").Value = 2297.25
$ws.Range("J137").Value = 5186
$ws.Range("K137").Value = 6891.75
$ws.Range("L137").Value = 15558
$ws.Range("M137").Value = -4341.75
$ws.Range("N137").Value = -20658

$ws.Range("H138").Value = 4990.6343
$ws.Range("J138").Value = 5603.515
$ws.Range("L138").Value = 16810.545
$ws.Range("N138").Value = -27090.545

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 949.3333
$ws.Range("I5").Value = 99.666664
$ws.Range("J5").Value = 1799
$ws.Range("K5").Value = 99.666664
$ws.Range("L5").Value = 1799
$ws.Range("M5").Value = 12.333336
$ws.Range("N5").Value = -2023

$ws.Range("H32").Value = 2119.8525
$ws.Range("I32").Value = 2148.6
$ws.Range("K32").Value = 2148.6
$ws.Range("M32").Value = -1861.6

$ws.Range("H63").Value = 10000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 10000
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -11372

$ws.Range("H66").Value = 10000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 50000
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -56864

$ws.Range("H132").Value = 10736.737
$ws.Range("I132").Value = 3999.8
$ws.Range("K132").Value = 11999.4
$ws.Range("M132").Value = -9469.400000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 949.3333
$ws.Range("I4").Value = 99.666664
$ws.Range("J4").Value = 1799
$ws.Range("K4").Value = 99.666664
$ws.Range("L4").Value = 1799
$ws.Range("M4").Value = 15.333336
$ws.Range("N4").Value = -2029

$ws.Range("H105").Value = 2398.0715
$ws.Range("I105").Value = 1000
$ws.Range("J105").Value = 2505.6155
$ws.Range("K105").Value = 1000
$ws.Range("L105").Value = 2505.6155
$ws.Range("M105").Value = 747
$ws.Range("N105").Value = -5999.6155

$ws.Range("H107").Value = 2040.0526
$ws.Range("I107").Value = 2007.7142
$ws.Range("K107").Value = 2007.7142
$ws.Range("M107").Value = -87.71419999999989

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1168.3684
$ws.Range("I94").Value = 878.7143
$ws.Range("J94").Value = 1337.3334
$ws.Range("K94").Value = 878.7143
$ws.Range("L94").Value = 1337.3334
$ws.Range("M94").Value = -427.7143
$ws.Range("N94").Value = -2239.3334

$ws.Range("H105").Value = 1821.7778
$ws.Range("I105").Value = 1780.1428
$ws.Range("K105").Value = 1780.1428
$ws.Range("M105").Value = -33.14280000000008

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1536.875
$ws.Range("I5").Value = 327.85715
$ws.Range("K5").Value = 983.5714499999999
$ws.Range("M5").Value = -871.5714499999999

$ws.Range("H22").Value = 1001
$ws.Range("I22").Value = 1001
$ws.Range("K22").Value = 3003
$ws.Range("M22").Value = -2834

$ws.Range("H27").Value = 1001
$ws.Range("I27").Value = 1001
$ws.Range("K27").Value = 3003
$ws.Range("M27").Value = -2901

$ws.Range("H121").Value = 911636.8
$ws.Range("I121").Value = 493.5
$ws.Range("J121").Value = 1432290.1
$ws.Range("K121").Value = 1480.5
$ws.Range("L121").Value = 4296870.300000001
$ws.Range("M121").Value = -170.5
$ws.Range("N121").Value = -4299490.300000001

$ws.Range("H132").Value = 3342.7144
$ws.Range("J132").Value = 3733.1667
$ws.Range("L132").Value = 33598.5003
$ws.Range("N132").Value = -38658.5003

$ws.Range("H135").Value = 1536.875
$ws.Range("I135").Value = 327.85715
$ws.Range("K135").Value = 2950.71435
$ws.Range("M135").Value = -415.7143499999997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 11372364
$ws.Range("I11").Value = 13754250
$ws.Range("K11").Value = 13754250
$ws.Range("M11").Value = -13754111

$ws.Range("H70").Value = 5717.9
$ws.Range("I70").Value = 5194.75
$ws.Range("J70").Value = 6066.6665
$ws.Range("K70").Value = 5194.75
$ws.Range("L70").Value = 6066.6665
$ws.Range("M70").Value = -4924.75
$ws.Range("N70").Value = -6606.6665

$ws.Range("H73").Value = 5717.9
$ws.Range("I73").Value = 5194.75
$ws.Range("J73").Value = 6066.6665
$ws.Range("K73").Value = 5194.75
$ws.Range("L73").Value = 6066.6665
$ws.Range("M73").Value = -4258.75
$ws.Range("N73").Value = -7938.6665

$ws.Range("H102").Value = 5232.7144
$ws.Range("I102").Value = 5635.269
$ws.Range("K102").Value = 5635.269
$ws.Range("M102").Value = -4013.269

$ws.Range("H126").Value = 2949.647
$ws.Range("I126").Value = 1678.5714
$ws.Range("J126").Value = 3839.4
$ws.Range("K126").Value = 5035.7142
$ws.Range("L126").Value = 11518.2
$ws.Range("M126").Value = -2565.7142
$ws.Range("N126").Value = -16458.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3811.8845
$ws.Range("I46").Value = 2682.647
$ws.Range("K46").Value = 2682.647
$ws.Range("M46").Value = -2494.647

$ws.Range("H122").Value = 622039.4399999999
$ws.Range("I122").Value = 479406.66
$ws.Range("K122").Value = 1438219.98
$ws.Range("M122").Value = -1435769.98

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 34486610
$ws.Range("I122").Value = 47621508
$ws.Range("J122").Value = 7500
$ws.Range("K122").Value = 142864524
$ws.Range("L122").Value = 22500
$ws.Range("M122").Value = -142862074
$ws.Range("N122").Value = -27400

